# RPA datasets push 2024-08-01
# Insert a new IPO record as the new top row (row 2), pushing all existing
# data rows down by one. The new record represents the most recent
# subscription date (2024-07-22) among the tracked IPOs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (first data row); this shifts
# every existing data row (2-24) down to (3-25) and keeps their values intact.
$ws.Rows.Item(2).Insert()
# Excel's row insert copies the format of the row above (the bold/bordered
# header) onto the new row -- clear it so the new data row matches the
# unstyled look of every other data row.
$ws.Rows.Item(2).ClearFormats()

$newRow = @(
    "2024-07-22",
    "피앤에스미캐닉스",
    "키움",
    "2024-07-25",
    "2024-07-31",
    29700000,
    1350000,
    "-",
    14000,
    17000,
    "-",
    22000,
    "-",
    "-",
    0,
    "-",
    "-",
    "1585.92 : 1",
    "-",
    "-"
)

# The date-like columns (A, D, E) must stay stored as plain text, matching
# the rest of the sheet -- force text format first so Excel doesn't coerce
# them into date serials.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

for ($i = 0; $i -lt $newRow.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $newRow[$i]
}
